$wb = $excel.ActiveWorkbook

# --- Fase4: replace "Cirio "Passata Verace" [AMC]" with "La Torrente, passata tradizionale [AMC]" and bump unit price 1.29 -> 1.39 ---
$wsFase4 = $wb.Worksheets.Item("Fase4")
$wsFase4.Range("A7").Value = "La Torrente, passata tradizionale [AMC]"
$wsFase4.Range("D7").Formula = "=1.39/0.7*B7/1000"
$wsFase4.Range("A8").Select()

# --- Fase5: replace "Cirio "Passata Verace" [AMC]" with "La Rosina, la passata di pomodoro [AMC]" and drop unit price 1.29 -> 0.85 ---
$wsFase5 = $wb.Worksheets.Item("Fase5")
$wsFase5.Range("A6").Value = "La Rosina, la passata di pomodoro [AMC]"
$wsFase5.Range("D6").Formula = "=0.85/0.7*B6/1000"
$wsFase5.Range("A6").Select()
